$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 255, shifting existing rows 255-264 down to 256-265.
$ws.Rows.Item(255).Insert()

# Populate the newly inserted row 255 with the new record.
$ws.Range("A255").Value = 4
$ws.Range("B255").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C255").Value = "Los Lagos"
$ws.Range("D255").Value = 44747
$ws.Range("D255").NumberFormat = $ws.Range("D256").NumberFormat
$ws.Range("E255").Value = 10
$ws.Range("F255").Value = "Fruta"
$ws.Range("G255").Value = 100108
$ws.Range("H255").Value = "Tropicales y subtropicales"
$ws.Range("I255").Value = 100108005
$ws.Range("J255").Value = "Piña"
$ws.Range("K255").Value = "Caramelo"
$ws.Range("L255").Value = "Tercera"
$ws.Range("M255").Value = 160
$ws.Range("N255").Value = 22000
$ws.Range("O255").Value = 22000
$ws.Range("P255").Value = 22000
$ws.Range("Q255").Value = "$/caja 16 unidades"
$ws.Range("R255").Value = "Ecuador"
$ws.Range("S255").Value = 1375
$ws.Range("T255").Value = 16
